$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.827802181243896
$ws.Range("B1").Value = 6.172130107879639
$ws.Range("C1").Value = 2.257118463516235
$ws.Range("D1").Value = 1.451746344566345
$ws.Range("E1").Value = 1.184410572052002
